# The workbook tracks weekly "Ajo" (garlic) price quotes for "Macroferia
# Regional de Talca". A new record (2021-09-21, serial 44460) needs to be
# inserted at row 112, which means every existing record currently on rows
# 112-161 shifts down by one row, and the freed-up row 112 receives the
# brand-new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 112; this pushes the old
# rows 112..161 down to 113..162 (and extends the sheet dimension to R162).
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly quote.
$ws.Cells.Item(112, 1).Value  = 5
$ws.Cells.Item(112, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value  = "Maule"
$ws.Cells.Item(112, 4).Value  = 44460
$ws.Cells.Item(112, 5).Value  = 7
$ws.Cells.Item(112, 6).Value  = 100112003
$ws.Cells.Item(112, 7).Value  = "Ajo"
$ws.Cells.Item(112, 8).Value  = "Chino"
$ws.Cells.Item(112, 9).Value  = "Primera"
$ws.Cells.Item(112, 10).Value = 300
$ws.Cells.Item(112, 11).Value = 15000
$ws.Cells.Item(112, 12).Value = 15000
$ws.Cells.Item(112, 13).Value = 15000
$ws.Cells.Item(112, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(112, 15).Value = "China"
$ws.Cells.Item(112, 16).Value = 1500
$ws.Cells.Item(112, 17).Value = 10
$ws.Cells.Item(112, 18).Value = "Hortaliza"
